$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 278, shifting the
# existing rows 278:344 down to 280:346 (dimension grows from T344 to T346).
$ws.Rows("278:279").Insert()

# Populate the two newly inserted rows (278 and 279) with the new data.
$ws.Range("A278").Value = 5
$ws.Range("B278").Value = "Macroferia Regional de Talca"
$ws.Range("C278").Value = "Maule"
$ws.Range("D278").Value = 44511
$ws.Range("E278").Value = 7
$ws.Range("F278").Value = "Fruta"
$ws.Range("G278").Value = 100102
$ws.Range("H278").Value = "Cítricos"
$ws.Range("I278").Value = 100102005
$ws.Range("J278").Value = "Naranja"
$ws.Range("K278").Value = "Navel Late"
$ws.Range("L278").Value = "Primera"
$ws.Range("M278").Value = 300
$ws.Range("N278").Value = 8000
$ws.Range("O278").Value = 8000
$ws.Range("P278").Value = 8000
$ws.Range("Q278").Value = "$/bandeja 15 kilos granel"
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value = 533
$ws.Range("T278").Value = 15

$ws.Range("A279").Value = 5
$ws.Range("B279").Value = "Macroferia Regional de Talca"
$ws.Range("C279").Value = "Maule"
$ws.Range("D279").Value = 44511
$ws.Range("E279").Value = 7
$ws.Range("F279").Value = "Fruta"
$ws.Range("G279").Value = 100102
$ws.Range("H279").Value = "Cítricos"
$ws.Range("I279").Value = 100102005
$ws.Range("J279").Value = "Naranja"
$ws.Range("K279").Value = "Olinda Valencia"
$ws.Range("L279").Value = "Primera"
$ws.Range("M279").Value = 250
$ws.Range("N279").Value = 8000
$ws.Range("O279").Value = 8000
$ws.Range("P279").Value = 8000
$ws.Range("Q279").Value = "$/bandeja 15 kilos granel"
$ws.Range("R279").Value = "Región de O'Higgins"
$ws.Range("S279").Value = 533
$ws.Range("T279").Value = 15
